# Auto-generated script to update Leve market-price columns (H-N) across all 8 sheets
# matching a scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC, row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 1049.25
$ws.Range("I2").Value = 949
$ws.Range("K2").Value = 949
$ws.Range("M2").Value = -836

$ws = $wb.Worksheets.Item(1)  # ALC, row 18 (Leve Item ID 5471)
$ws.Range("H18").Value = 406.8889
$ws.Range("I18").Value = 406.8889
$ws.Range("K18").Value = 406.8889
$ws.Range("M18").Value = -122.8889

$ws = $wb.Worksheets.Item(1)  # ALC, row 32 (Leve Item ID 5484)
$ws.Range("H32").Value = 7315.5
$ws.Range("J32").Value = 7520.5
$ws.Range("L32").Value = 7520.5
$ws.Range("N32").Value = -8172.5

$ws = $wb.Worksheets.Item(1)  # ALC, row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 400.7143
$ws.Range("I33").Value = 501.3
$ws.Range("K33").Value = 501.3
$ws.Range("M33").Value = -272.3

$ws = $wb.Worksheets.Item(1)  # ALC, row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 2566.3333
$ws.Range("I43").Value = 2350
$ws.Range("K43").Value = 2350
$ws.Range("M43").Value = -2281

$ws = $wb.Worksheets.Item(1)  # ALC, row 86 (Leve Item ID 12603)
$ws.Range("H86").Value = 1582.8334
$ws.Range("J86").Value = 1499.75
$ws.Range("L86").Value = 1499.75
$ws.Range("N86").Value = -3745.75

$ws = $wb.Worksheets.Item(1)  # ALC, row 89 (Leve Item ID 12603)
$ws.Range("H89").Value = 1582.8334
$ws.Range("J89").Value = 1499.75
$ws.Range("L89").Value = 7498.75
$ws.Range("N89").Value = -18730.75

$ws = $wb.Worksheets.Item(1)  # ALC, row 92 (Leve Item ID 19901)
$ws.Range("H92").Value = 593.2727
$ws.Range("I92").Value = 642.1053000000001
$ws.Range("J92").Value = 284
$ws.Range("K92").Value = 642.1053000000001
$ws.Range("L92").Value = 284
$ws.Range("M92").Value = 605.8946999999999
$ws.Range("N92").Value = -2780

$ws = $wb.Worksheets.Item(1)  # ALC, row 106 (Leve Item ID 19903)
$ws.Range("H106").Value = 2931.6667
$ws.Range("I106").Value = 2897.5
$ws.Range("K106").Value = 2897.5
$ws.Range("M106").Value = -2266.5

$ws = $wb.Worksheets.Item(1)  # ALC, row 112 (Leve Item ID 27960)
$ws.Range("H112").Value = 41669588
$ws.Range("J112").Value = 71431930
$ws.Range("L112").Value = 214295790
$ws.Range("N112").Value = -214298006

$ws = $wb.Worksheets.Item(1)  # ALC, row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 23258044
$ws.Range("I138").Value = 884.0526
$ws.Range("J138").Value = 41669960
$ws.Range("K138").Value = 2652.1578
$ws.Range("L138").Value = 125009880
$ws.Range("M138").Value = 2487.8422
$ws.Range("N138").Value = -125020160

$ws = $wb.Worksheets.Item(1)  # ALC, row 141 (Leve Item ID 44161)
$ws.Range("H141").Value = 1861
$ws.Range("I141").Value = 1047.5
$ws.Range("J141").Value = 2403.3333
$ws.Range("K141").Value = 3142.5
$ws.Range("L141").Value = 7209.999899999999
$ws.Range("M141").Value = 2037.5
$ws.Range("N141").Value = -17569.9999

$ws = $wb.Worksheets.Item(2)  # ARM, row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 3177.35
$ws.Range("I2").Value = 2596.5625
$ws.Range("J2").Value = 5500.5
$ws.Range("K2").Value = 2596.5625
$ws.Range("L2").Value = 5500.5
$ws.Range("M2").Value = -2483.5625
$ws.Range("N2").Value = -5726.5

$ws = $wb.Worksheets.Item(2)  # ARM, row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 5347.8945
$ws.Range("I32").Value = 5741.697
$ws.Range("K32").Value = 5741.697
$ws.Range("M32").Value = -5454.697

$ws = $wb.Worksheets.Item(2)  # ARM, row 41 (Leve Item ID 2501)
$ws.Range("H41").Value = 4937.4287
$ws.Range("I41").Value = 4083.3333
$ws.Range("J41").Value = 10062
$ws.Range("K41").Value = 4083.3333
$ws.Range("L41").Value = 10062
$ws.Range("M41").Value = -3669.3333
$ws.Range("N41").Value = -10890

$ws = $wb.Worksheets.Item(2)  # ARM, row 60 (Leve Item ID 3883)
$ws.Range("H60").Value = 39334.332
$ws.Range("I60").Value = 16499.5
$ws.Range("J60").Value = 85004
$ws.Range("K60").Value = 16499.5
$ws.Range("L60").Value = 85004
$ws.Range("M60").Value = -15766.5
$ws.Range("N60").Value = -86470

$ws = $wb.Worksheets.Item(2)  # ARM, row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 29897.584
$ws.Range("I110").Value = 32559.092
$ws.Range("J110").Value = 621
$ws.Range("K110").Value = 32559.092
$ws.Range("L110").Value = 621
$ws.Range("M110").Value = -30514.092
$ws.Range("N110").Value = -4711

$ws = $wb.Worksheets.Item(2)  # ARM, row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 3177.35
$ws.Range("I116").Value = 2596.5625
$ws.Range("J116").Value = 5500.5
$ws.Range("K116").Value = 2596.5625
$ws.Range("L116").Value = 5500.5
$ws.Range("M116").Value = -302.5625
$ws.Range("N116").Value = -10088.5

$ws = $wb.Worksheets.Item(2)  # ARM, row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 119783.05
$ws.Range("I132").Value = 6163.5557
$ws.Range("K132").Value = 18490.6671
$ws.Range("M132").Value = -15960.6671

$ws = $wb.Worksheets.Item(3)  # BSM, row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 3177.35
$ws.Range("I3").Value = 2596.5625
$ws.Range("J3").Value = 5500.5
$ws.Range("K3").Value = 2596.5625
$ws.Range("L3").Value = 5500.5
$ws.Range("M3").Value = -2482.5625
$ws.Range("N3").Value = -5728.5

$ws = $wb.Worksheets.Item(4)  # CRP, row 7 (Leve Item ID 5361)
$ws.Range("H7").Value = 147.90909
$ws.Range("I7").Value = 43.6
$ws.Range("K7").Value = 43.6
$ws.Range("M7").Value = 69.40000000000001

$ws = $wb.Worksheets.Item(4)  # CRP, row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 2907.1072
$ws.Range("I31").Value = 1842.3077
$ws.Range("J31").Value = 3829.9333
$ws.Range("K31").Value = 1842.3077
$ws.Range("L31").Value = 3829.9333
$ws.Range("M31").Value = -1547.3077
$ws.Range("N31").Value = -4419.933300000001

$ws = $wb.Worksheets.Item(4)  # CRP, row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2907.1072
$ws.Range("I34").Value = 1842.3077
$ws.Range("J34").Value = 3829.9333
$ws.Range("K34").Value = 1842.3077
$ws.Range("L34").Value = 3829.9333
$ws.Range("M34").Value = -1640.3077
$ws.Range("N34").Value = -4233.933300000001

$ws = $wb.Worksheets.Item(4)  # CRP, row 62 (Leve Item ID 12580)
$ws.Range("H62").Value = 3827.5715
$ws.Range("I62").Value = 3746.6667
$ws.Range("K62").Value = 3746.6667
$ws.Range("M62").Value = -3122.6667

$ws = $wb.Worksheets.Item(4)  # CRP, row 65 (Leve Item ID 12580)
$ws.Range("H65").Value = 3827.5715
$ws.Range("I65").Value = 3746.6667
$ws.Range("K65").Value = 18733.3335
$ws.Range("M65").Value = -15613.3335

$ws = $wb.Worksheets.Item(4)  # CRP, row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 3827.3044
$ws.Range("I132").Value = 4006.6191
$ws.Range("J132").Value = 1944.5
$ws.Range("K132").Value = 12019.8573
$ws.Range("L132").Value = 5833.5
$ws.Range("M132").Value = -9489.8573
$ws.Range("N132").Value = -10893.5

$ws = $wb.Worksheets.Item(4)  # CRP, row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 1498.75
$ws.Range("I134").Value = 1192.25
$ws.Range("K134").Value = 3576.75
$ws.Range("M134").Value = -1041.75

$ws = $wb.Worksheets.Item(5)  # CUL, row 4 (Leve Item ID 4650)
$ws.Range("H4").Value = 12447430
$ws.Range("I4").Value = 4439754.5
$ws.Range("K4").Value = 13319263.5
$ws.Range("M4").Value = -13319151.5

$ws = $wb.Worksheets.Item(5)  # CUL, row 11 (Leve Item ID 4745)
$ws.Range("H11").Value = 60.42857
$ws.Range("I11").Value = 70
$ws.Range("J11").Value = 3
$ws.Range("K11").Value = 210
$ws.Range("L11").Value = 9
$ws.Range("M11").Value = -70
$ws.Range("N11").Value = -289

$ws = $wb.Worksheets.Item(5)  # CUL, row 121 (Leve Item ID 27878)
$ws.Range("H121").Value = 2000
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 2000
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 6000
$ws.Range("N121").Value = -8620
$ws.Range("M121").ClearContents()

$ws = $wb.Worksheets.Item(6)  # GSM, row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 4326.7856
$ws.Range("I102").Value = 2145.7368
$ws.Range("K102").Value = 2145.7368
$ws.Range("M102").Value = -523.7368000000001

$ws = $wb.Worksheets.Item(6)  # GSM, row 125 (Leve Item ID 34291)
$ws.Range("H125").Value = 55000
$ws.Range("J125").Value = 55000
$ws.Range("L125").Value = 55000
$ws.Range("N125").Value = -59920

$ws = $wb.Worksheets.Item(6)  # GSM, row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 3371.375
$ws.Range("I126").Value = 2831.3333
$ws.Range("K126").Value = 8493.999899999999
$ws.Range("M126").Value = -6023.999899999999

$ws = $wb.Worksheets.Item(7)  # LTW, row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 2631.4707
$ws.Range("I22").Value = 2019.8
$ws.Range("K22").Value = 2019.8
$ws.Range("M22").Value = -1724.8

$ws = $wb.Worksheets.Item(7)  # LTW, row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 2631.4707
$ws.Range("I27").Value = 2019.8
$ws.Range("K27").Value = 2019.8
$ws.Range("M27").Value = -1912.8

$ws = $wb.Worksheets.Item(7)  # LTW, row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 5577.6
$ws.Range("I40").Value = 4472
$ws.Range("K40").Value = 4472
$ws.Range("M40").Value = -4336

$ws = $wb.Worksheets.Item(7)  # LTW, row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 2875.92

$ws = $wb.Worksheets.Item(7)  # LTW, row 82 (Leve Item ID 12565)
$ws.Range("H82").Value = 3728.3333
$ws.Range("I82").Value = 3417.75
$ws.Range("K82").Value = 3417.75
$ws.Range("M82").Value = -3056.75

$ws = $wb.Worksheets.Item(7)  # LTW, row 85 (Leve Item ID 12565)
$ws.Range("H85").Value = 3728.3333
$ws.Range("I85").Value = 3417.75
$ws.Range("K85").Value = 3417.75
$ws.Range("M85").Value = -2169.75

$ws = $wb.Worksheets.Item(7)  # LTW, row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 3285.2
$ws.Range("I93").Value = 2725
$ws.Range("K93").Value = 2725
$ws.Range("M93").Value = -1477

$ws = $wb.Worksheets.Item(7)  # LTW, row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 1837.238
$ws.Range("I132").Value = 1343.2222
$ws.Range("K132").Value = 4029.6666
$ws.Range("M132").Value = -1499.6666

$ws = $wb.Worksheets.Item(7)  # LTW, row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 4500.64
$ws.Range("I136").Value = 3022.4348
$ws.Range("J136").Value = 21500
$ws.Range("K136").Value = 9067.304400000001
$ws.Range("L136").Value = 64500
$ws.Range("M136").Value = -6517.304400000001
$ws.Range("N136").Value = -69600

$ws = $wb.Worksheets.Item(8)  # WVR, row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 3650.3333
$ws.Range("I122").Value = 2698.2856
$ws.Range("K122").Value = 8094.8568
$ws.Range("M122").Value = -5644.8568

$ws = $wb.Worksheets.Item(8)  # WVR, row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 1185.2941
$ws.Range("I132").Value = 1160.6061
$ws.Range("K132").Value = 3481.8183
$ws.Range("M132").Value = -951.8182999999999

Write-Output "Updated market-price columns across all sheets."
